# Apply German (de) field translations and related fixes to the fields.xlsx workbook
# per commit "Added german fields translation".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: gc_fields_display  (columns: A=field, B=en, C=fr, D=de, E=it, F=es, G=pt, H=ja, I=zh)
# ---------------------------------------------------------------------------
$wsDisplay = $wb.Worksheets.Item("gc_fields_display")

# New German translations (cells that previously had no D value)
$wsDisplay.Range("D2").Value  = "Beginn Breitengrad"
$wsDisplay.Range("D3").Value  = "Beginn Längengrad"
$wsDisplay.Range("D7").Value  = "Beginn Zeitstempel"
$wsDisplay.Range("D12").Value = "Ende Breitengrad"
$wsDisplay.Range("D13").Value = "Ende Längengrad"
$wsDisplay.Range("D17").Value = "Ende Zeitstempel"
$wsDisplay.Range("D31").Value = "Max. korrigierte Höhe"
$wsDisplay.Range("D42").Value = "Max. Leistung"
$wsDisplay.Range("D43").Value = "Max. Leistung 20min"
$wsDisplay.Range("D44").Value = "Max. relative Laufökonomie"
$wsDisplay.Range("D51").Value = "Max. vertikales Verhältnis "
$wsDisplay.Range("D52").Value = "Max. vertikale Geschw."
$wsDisplay.Range("D55").Value = "Min. Trittfrequenz "
$wsDisplay.Range("D56").Value = "Min. korrigierter Anstieg"
$wsDisplay.Range("D62").Value = "Min. Puls"
$wsDisplay.Range("D65").Value = "Min. Tempo"
$wsDisplay.Range("D66").Value = "Min. Leistung"
$wsDisplay.Range("D68").Value = "Min. Schrittfrequenz"
$wsDisplay.Range("D86").Value = "Aktivitätskalorien"
$wsDisplay.Range("D89").Value = "Intensitätsfaktor"
$wsDisplay.Range("D112").Value = "TSS"
$wsDisplay.Range("D113").Value = "FTP"

# Updated German translations (replacing existing D values)
$wsDisplay.Range("D19").Value  = "Höhengewinn"
$wsDisplay.Range("D26").Value  = "Höhenverlust"
$wsDisplay.Range("D30").Value  = "Max. Trittfrequenz (Rad)"
$wsDisplay.Range("D32").Value  = "Max. Schrittfrequenz"
$wsDisplay.Range("D35").Value  = "Max. Schrittfrequenz"
$wsDisplay.Range("D38").Value  = "Max. Puls"
$wsDisplay.Range("D41").Value  = "Max. Tempo"
$wsDisplay.Range("D45").Value  = "Max. Schrittfrequenz"
$wsDisplay.Range("D48").Value  = "Max. Zugfrequenz"
$wsDisplay.Range("D69").Value  = "Min. Geschw."
$wsDisplay.Range("D70").Value  = "Min. Züge"
$wsDisplay.Range("D84").Value  = "Absolvierte Zeit"
$wsDisplay.Range("D111").Value = "Trainingseffekt"
$wsDisplay.Range("D116").Value = "Ø Trittfrequenz (Rad)"
$wsDisplay.Range("D117").Value = "Ø Schrittfrequenz"
$wsDisplay.Range("D118").Value = "Ø Effizienz"
$wsDisplay.Range("D121").Value = "Ø Schrittfrequenz"
$wsDisplay.Range("D123").Value = "Ø Bodenkontaktzeit"
$wsDisplay.Range("D124").Value = "Ø Puls"
$wsDisplay.Range("D130").Value = "Ø Bewegungstempo"
$wsDisplay.Range("D131").Value = "Ø Bewegungsgeschw."
$wsDisplay.Range("D133").Value = "Ø Tempo"
$wsDisplay.Range("D137").Value = "Ø Schrittfrequenz"
$wsDisplay.Range("D138").Value = "Ø Geschw."
$wsDisplay.Range("D141").Value = "Ø Zuglänge"
$wsDisplay.Range("D144").Value = "Ø Anzahl von Zügen"
$wsDisplay.Range("D145").Value = "Ø Zugfrequenz"
$wsDisplay.Range("D146").Value = "Ø SWOLF-Wert"
$wsDisplay.Range("D147").Value = "Ø vertikale Bewegung"

# ---------------------------------------------------------------------------
# Sheet: gc_fields_order  (SumStrokes rows: fix swim/dynamics grouping)
# ---------------------------------------------------------------------------
$wsOrder = $wb.Worksheets.Item("gc_fields_order")

$wsOrder.Range("B92").Value = "swimming"
$wsOrder.Range("C93").Value = "dynamics"
$wsOrder.Range("D93").Value = 2
$wsOrder.Range("B93").ClearContents()

# ---------------------------------------------------------------------------
# Sheet: gc_activity_types  (columns include G=de)
# ---------------------------------------------------------------------------
$wsActivity = $wb.Worksheets.Item("gc_activity_types")

$wsActivity.Range("G7").Value  = "Virtuelles Laufen"
$wsActivity.Range("G8").Value  = "Hindernislauf"
$wsActivity.Range("G9").Value  = "Indoor-Laufen"
$wsActivity.Range("G19").Value = "Virtuelles Radfahren"
$wsActivity.Range("G41").Value = "Tennis"
$wsActivity.Range("G42").Value = "Treppensteigen"
$wsActivity.Range("G43").Value = "Stoppuhr"
$wsActivity.Range("G44").Value = "Autorennen"
$wsActivity.Range("G45").Value = "Atemübung"
$wsActivity.Range("G58").Value = "Multi-Sport"
$wsActivity.Range("G59").Value = "Schritte"
$wsActivity.Range("G60").Value = "Tauchen"
$wsActivity.Range("G61").Value = "Sicherheit"
$wsActivity.Range("G62").Value = "Wintersport"
$wsActivity.Range("G70").Value = "Pilates"
$wsActivity.Range("G71").Value = "Yoga"
$wsActivity.Range("G77").Value = "Einzelgespräche-Tauchgang"
$wsActivity.Range("G78").Value = "Multigas-Tauchgang"
$wsActivity.Range("G79").Value = "Tiefenmesser-Tauchgang"
$wsActivity.Range("G80").Value = "Apnoetauchgang"
$wsActivity.Range("G81").Value = "Apnoejagd"
$wsActivity.Range("G82").Value = "CCR-Tauchgang"
$wsActivity.Range("G83").Value = "Hilfe"
$wsActivity.Range("G84").Value = "Unfall-Benachrichtigungen"
